$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new issue row (row 7)
$ws.Range("A7").Value = "Bot stopped responding at unit selection page"
$ws.Range("B7").Value = "Bugfix"
$ws.Range("C7").Value = "Yes"
$ws.Range("D7").Value = "Bot just stops at the unit selection page. "

$ws.Range("A7").HorizontalAlignment = -4131
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("C7").HorizontalAlignment = -4108
$ws.Range("D7").HorizontalAlignment = -4131

# Update the view's top left cell (B1) and active selection (C7)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C7").Select()
